$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("F1").Value = "SAT Math"
$ws.Range("G1").Value = "SAT Reading"

# New header cells H1/I1 need the same bold/border header style as the
# existing header cells, so copy formatting from G1 before writing values.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "SAT Combined"
$ws.Range("I1").Value = "GPA"

# Update row 2 (existing data row)
$ws.Range("A2").Value = "Test One"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 8000
$ws.Range("D2").Value = "All"
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = 600
$ws.Range("G2").Value = 400
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 4

# Add new row 3
$ws.Range("A3").Value = "Cool Kids Club"
$ws.Range("B3").Value = "'1"
$ws.Range("C3").Value = "'11111111"
$ws.Range("D3").Value = "Computer Science and Engineering"
$ws.Range("E3").Value = 36
$ws.Range("F3").Value = 800
$ws.Range("G3").Value = 800
$ws.Range("H3").Value = 1600
$ws.Range("I3").Value = 5
